$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "holly added S.GISH to harvester in bioSamples" -- the harvester column (B)
# for every data row (2-19) changes from "Retrofitted_1045" to "S.GISH".
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 2).Value = "S.GISH"
}

# Select column B, matching the selection left behind in the saved file.
$ws.Range("B:B").Select()
